$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metrics value set (columns B..Q) that every data row (2-26) now shares
$newValues = @(
    [double]"0.999999900753416",        # B r2
    [double]"0.6876245000598846",        # C r2_sup
    [double]"0.9999997502514453",        # D r2_test
    [double]"0.9999999448432606",        # E r2_val
    [double]"0.9999999426895874",        # F r2_vt
    [double]"5.891705945054375e-08",     # G mse
    [double]"0.1854395906834327",        # H mse_sup
    [double]"8.363343981512521e-08",     # I mse_test
    [double]"1.444494337408425e-08",     # J mse_val
    [double]"4.900472405255845e-08",     # K mse_vt
    [double]"9.931498004524596e-05",     # L mape
    [double]"0.0002427283655664161",     # M rmse
    [double]"1.000000140112824",         # N r2_adj
    [double]"0.0002530618146055635",     # O rsd
    [double]"115.2942703081797",         # P aic
    [double]"165.268179127776"           # Q bic
)

# Apply the new value set to every data row (rows 2 through 26), columns B..Q
for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}

# Update the model labels in column A: a 3-way cyclic rename
# row 2:  model_24_6_12 -> model_24_6_0
# row 14: model_24_6_24 -> model_24_6_12
# row 26: model_24_6_0  -> model_24_6_24
$ws.Range("A2").Value = "model_24_6_0"
$ws.Range("A14").Value = "model_24_6_12"
$ws.Range("A26").Value = "model_24_6_24"
